$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single Column Numbers")

# The new (type-aware, stable) sort implementation changes the relative
# order of the two rows that tie on the "Numbers" column (A2:A3 = 5, 5):
# the row that used to be row 2 ("Jacques", visible) and the row that
# used to be row 3 ("Alex", hidden by the autofilter) trade places.
#
# Unhide the destination row first, so the "Value" write below does not
# happen on a still-hidden row (which would otherwise trigger a spurious
# autofit row-height recalculation on save).
$ws.Rows(3).Hidden = $false

# Swap the "Names" values held by rows 2 and 3.
$ws.Range("B2").Value = "Alex"
$ws.Range("B3").Value = "Jacques"

# Row 2 ("Alex") is now the one filtered out by the autofilter, so hide
# it last (after its value has already been updated).
$ws.Rows(2).Hidden = $true
